# Add a new "2022-Q3" worksheet (by cloning the "2021-Q3" sheet, which already
# has the right 7-fund / 8-row layout + formatting) and fill it with the new
# quarter's fund-holding data, then insert the corresponding summary row into
# the "总计" (Total) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by copying "2021-Q3" (currently item 5)
#    so it inherits identical styles/borders/column layout, and place it
#    right before the current "2022-Q2" sheet (item 2).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(5)
$beforeSheet = $wb.Worksheets.Item(2)
$template.Copy($beforeSheet)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Header row
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Data rows (A=index, B=code, C=name, D=scale, E=stock position, F=position pct, G=held value, H=rank)
$q3Rows = @(
    @(0, "160613", "鹏华盛世创新混合（LOF）", "2.53", "92.08", "4.86", "0.1230", 2),
    @(1, "012640", "鹏华稳健鸿利一年持有期混合A", "2.61", "92.98", "4.35", "0.1135", 6),
    @(2, "008134", "鹏华优选价值股票", "1.80", "92.72", "4.87", "0.0877", 5),
    @(3, "011574", "鹏华领航一年持有期混合A", "1.20", "92.84", "4.41", "0.0529", 6),
    @(4, "011575", "鹏华领航一年持有期混合C", "0.91", "92.84", "4.41", "0.0401", 6),
    @(5, "001219", "上投摩根动态多因子策略混合", "1.02", "90.78", "3.31", "0.0338", 8),
    @(6, "012641", "鹏华稳健鸿利一年持有期混合C", "0.10", "92.98", "4.35", "0.0044", 6)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (Total) summary sheet: insert the 2022-Q3 row at the
#    top of the data (row 2) and push the rest down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$totalRows = @(
    @("2022-Q3", 7, 0.46),
    @("2022-Q2", 2, 0.19),
    @("2022-Q1", 4, 0.38),
    @("2021-Q4", 4, 0.42),
    @("2021-Q3", 7, 1.45),
    @("2021-Q2", 5, 1.35),
    @("2021-Q1", 4, 0.53),
    @("2020-Q4", 1, 0.19)
)

$r = 2
$idx = 0
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $idx
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}
